$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Give row 10 (A10:F10) a closed bottom border, matching row 9's "full box" look ---
# (G10 already has the full box border.)
$row10 = $ws.Range("A10:F10")
$row10.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$row10.Borders.Item(9).Weight = 2      # xlThin

# --- New row 11: full thin box border around A11:G11 ---
$row11 = $ws.Range("A11:G11")
$row11.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$row11.Borders.Item(7).Weight = 2
$row11.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$row11.Borders.Item(8).Weight = 2
$row11.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$row11.Borders.Item(9).Weight = 2
$row11.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$row11.Borders.Item(10).Weight = 2
$row11.Borders.Item(11).LineStyle = 1  # xlInsideVertical
$row11.Borders.Item(11).Weight = 2

$row11.HorizontalAlignment = -4108     # xlCenter
$row11.VerticalAlignment = -4108       # xlCenter
$row11.WrapText = $false

$ws.Range("C11").WrapText = $true
$ws.Range("F11").WrapText = $true
$ws.Range("F11").NumberFormat = "000"

# --- Values for the new "Sacudir Elementos" fraction row ---
$ws.Range("A11").Value = "FR-SA-001"
$ws.Range("B11").Value = "Sacudir Elementos"
$ws.Range("C11").Value = "Revisar que todos los objetos queden correctamente acomodados."
$ws.Range("D11").Value = "FR"
$ws.Range("E11").Value = "SA"
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = "Sacudir Elementos"

$ws.Range("G14").Select()
